# Applies the updated classification values from the commit
# "Merging prediction and classification into a single repo"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - Control 33
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.410845134592065
$ws.Range("E3").Value = 0.410845134592065

# Row 4 - Control 36
$ws.Range("D4").Value = 0.9924099073764798
$ws.Range("E4").Value = 0.9924099073764798

# Row 5 - Control 49
$ws.Range("D5").Value = [double]"1.439065369307728E-09"
$ws.Range("E5").Value = [double]"1.439065369307728E-09"

# Row 6 - Control 2
$ws.Range("D6").Value = 0.001142764819202968
$ws.Range("E6").Value = 0.001142764819202968

# Row 8 - MDD 12
$ws.Range("D8").Value = 0.9999134664762226
$ws.Range("E8").Value = [double]"8.653352377741008E-05"

# Row 9 - MDD 53
$ws.Range("D9").Value = 0.999999999999996
$ws.Range("E9").Value = [double]"3.996802888650564E-15"

# Row 10 - MDD 29
$ws.Range("D10").Value = [double]"1.716369506482791E-09"
$ws.Range("E10").Value = 0.9999999982836305

# Row 11 - MDD 55
$ws.Range("D11").Value = [double]"2.432130054037753E-36"
$ws.Range("F11").Value = 14.68372821807861
$ws.Range("G11").Value = 0.6
